# Daily Conversation.docx edit:
#  - "M" + "any test failed" (two runs) -> single run "Many test failed"
#  - Append the rest of the Slack conversation as new paragraphs right
#    after that line (13 new paragraphs, some blank separators).

function Set-RunXml($range, $text) {
    # Replace $range's content with a single run containing $text, while
    # keeping the same boilerplate <w:rPr/> this document always uses.
    $escaped = $text -replace '&', '&amp;'
    $escaped = $escaped -replace '<', '&lt;'
    $escaped = $escaped -replace '>', '&gt;'
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
        + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
        + '<w:body><w:p><w:r><w:rPr/><w:t>' + $escaped + '</w:t></w:r></w:p></w:body></w:document>' `
        + '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

$d = $word.ActiveDocument

# Locate the "Many test failed" paragraph (currently stored as two runs:
# "M" and "any test failed"); search from the end since it's near the
# bottom of the document.
$paras = $d.Paragraphs
$count = $paras.Count
$targetIndex = -1
for ($i = $count; $i -ge 1; $i--) {
    $pp = $paras.Item($i)
    if ($pp.Range.Text.Trim() -eq "Many test failed") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Many test failed' paragraph"
}

# Merge the two runs "M" / "any test failed" into a single run, keeping
# the paragraph's trailing mark untouched.
$p = $paras.Item($targetIndex)
$mergeRange = $d.Range($p.Range.Start, $p.Range.End - 1)
Set-RunXml $mergeRange "Many test failed"

# New paragraphs to insert right after it (Slack conversation continues).
$newParaTexts = @(
    "",
    "Sakibur Rahaman 7:00 PM",
    "@here Are we going to have any pre-demo meetings today?",
    "",
    "Travis Emslander:spiral_calendar_pad: 7:10 PM",
    "@here Andrew seems to be running late. Can the team put together the items they'd like to talk about and/or show? Here's the full list: https://fieldnation.atlassian.net/wiki/spaces/DEV/pages/761856267/Sprint+19.21",
    "",
    "John Vogt:speech_balloon: 7:11 PM",
    "Ya, let's get into the standup and chat quick, can someone present the list Travis has above",
    "",
    "Andrew Kandels 7:31 PM",
    "@here sorry have a real sick kid which caused me to sleep through my alarm. here WFH though, ill present the recruitments stuff, looks like you guys divied up the rest",
    ""
)

$insertAfterIndex = $targetIndex
foreach ($t in $newParaTexts) {
    $paras = $d.Paragraphs
    $anchorPara = $paras.Item($insertAfterIndex)
    $null = $anchorPara.Range.InsertParagraphAfter()
    $insertAfterIndex = $insertAfterIndex + 1
    if ($t -ne "") {
        $paras = $d.Paragraphs
        $newPara = $paras.Item($insertAfterIndex)
        $newPara.Range.Text = $t
    }
}

Write-Output "Inserted $($newParaTexts.Count) paragraphs after paragraph $targetIndex; total paragraphs now $($d.Paragraphs.Count)"
